$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.022559500331671
$ws.Range("D2").Value = 1.028240752498946
$ws.Range("E2").Value = 1.023296285314085
$ws.Range("I2").Value = 1.031292256682527
$ws.Range("J2").Value = 1.02774465489041
$ws.Range("K2").Value = 1.03105814358181
$ws.Range("L2").Value = 1.026128140484702
$ws.Range("N2").Value = 1.013326145891841
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.023268414497513
$ws.Range("D3").Value = 1.028757867005622
$ws.Range("E3").Value = 1.023891552078117
$ws.Range("I3").Value = 1.031404072018624
$ws.Range("J3").Value = 1.028092919966117
$ws.Range("K3").Value = 1.031383937220203
$ws.Range("L3").Value = 1.026530833216134
$ws.Range("N3").Value = 1.013441430053179
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.023727788980767
$ws.Range("D4").Value = 1.029092973938511
$ws.Range("E4").Value = 1.02427769913853
$ws.Range("I4").Value = 1.031475514945559
$ws.Range("J4").Value = 1.028318231701515
$ws.Range("K4").Value = 1.031594532211163
$ws.Range("L4").Value = 1.026791674326228
$ws.Range("N4").Value = 1.013516002383445
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.023921065817073
$ws.Range("D5").Value = 1.029233970396204
$ws.Range("E5").Value = 1.024440265402363
$ws.Range("I5").Value = 1.031505331152814
$ws.Range("J5").Value = 1.028412942036504
$ws.Range("K5").Value = 1.031683013401435
$ws.Range("L5").Value = 1.026901395518175
$ws.Range("N5").Value = 1.013547346318703
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.023953526923395
$ws.Range("D6").Value = 1.029257651153379
$ws.Range("E6").Value = 1.024467574421888
$ws.Range("I6").Value = 1.031510324595928
$ws.Range("J6").Value = 1.028428843650122
$ws.Range("K6").Value = 1.031697866637035
$ws.Range("L6").Value = 1.02691982187361
$ws.Range("N6").Value = 1.01355260872181
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.02373037094576
$ws.Range("D7").Value = 1.029094857480925
$ws.Range("E7").Value = 1.024279870455654
$ws.Range("I7").Value = 1.031475914210296
$ws.Range("J7").Value = 1.028319497269484
$ws.Range("K7").Value = 1.031595714710961
$ws.Range("L7").Value = 1.026793140178224
$ws.Range("N7").Value = 1.013516421227882
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.022798943347535
$ws.Range("D8").Value = 1.028415409176045
$ws.Range("E8").Value = 1.023497256469095
$ws.Range("I8").Value = 1.031330232967939
$ws.Range("J8").Value = 1.027862359968063
$ws.Range("K8").Value = 1.031168290797446
$ws.Range("L8").Value = 1.026264175146489
$ws.Range("N8").Value = 1.013365111447912
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.021162794114935
$ws.Range("D9").Value = 1.027222047598256
$ws.Range("E9").Value = 1.02212570979588
$ws.Range("I9").Value = 1.031066595492564
$ws.Range("J9").Value = 1.027056586997696
$ws.Range("K9").Value = 1.030413528614407
$ws.Range("L9").Value = 1.025334221632805
$ws.Range("N9").Value = 1.013098319533667
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.020075613639945
$ws.Range("D10").Value = 1.026429225202065
$ws.Range("E10").Value = 1.021216526429121
$ws.Range("I10").Value = 1.030886224749195
$ws.Range("J10").Value = 1.026519326461464
$ws.Range("K10").Value = 1.029909371199445
$ws.Range("L10").Value = 1.024715789705058
$ws.Range("N10").Value = 1.012920375685183
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.019605729777308
$ws.Range("D11").Value = 1.026086602803255
$ws.Range("E11").Value = 1.02082409558772
$ws.Range("I11").Value = 1.030807037959869
$ws.Range("J11").Value = 1.02628668546959
$ws.Range("K11").Value = 1.029690850270821
$ws.Range("L11").Value = 1.024448385809254
$ws.Range("N11").Value = 1.01284331043334
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.019431327249483
$ws.Range("D12").Value = 1.025959441092668
$ws.Range("E12").Value = 1.020678519522587
$ws.Range("I12").Value = 1.030777462221318
$ws.Range("J12").Value = 1.026200273114445
$ws.Range("K12").Value = 1.029609650651503
$ws.Range("L12").Value = 1.024349118950473
$ws.Range("N12").Value = 1.012814683291907
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.019468731118367
$ws.Range("D13").Value = 1.025986712993796
$ws.Range("E13").Value = 1.020709737439665
$ws.Range("I13").Value = 1.030783813646776
$ws.Range("J13").Value = 1.026218808796028
$ws.Range("K13").Value = 1.02962706963484
$ws.Range("L13").Value = 1.024370409340013
$ws.Range("N13").Value = 1.012820823978944
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.019591310878859
$ws.Range("D14").Value = 1.02607608945301
$ws.Range("E14").Value = 1.020812058331643
$ws.Range("I14").Value = 1.030804596525213
$ws.Range("J14").Value = 1.026279542571493
$ws.Range("K14").Value = 1.029684138913388
$ws.Range("L14").Value = 1.024440179169998
$ws.Range("N14").Value = 1.01284094413627
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.01966685403852
$ws.Range("D15").Value = 1.026131171017345
$ws.Range("E15").Value = 1.020875126881912
$ws.Range("I15").Value = 1.030817380067605
$ws.Range("J15").Value = 1.026316962812224
$ws.Range("K15").Value = 1.029719297086257
$ws.Range("L15").Value = 1.024483174510956
$ws.Range("N15").Value = 1.012853340621752
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.020106816852462
$ws.Range("D16").Value = 1.026451978321934
$ws.Range("E16").Value = 1.021242597327571
$ws.Range("I16").Value = 1.030891457321049
$ws.Range("J16").Value = 1.026534766122063
$ws.Range("K16").Value = 1.02992386925651
$ws.Range("L16").Value = 1.024733544587713
$ws.Range("N16").Value = 1.01292548999049
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.020383029025281
$ws.Range("D17").Value = 1.026653394453047
$ws.Range("E17").Value = 1.021473438454866
$ws.Range("I17").Value = 1.030937634072149
$ws.Range("J17").Value = 1.026671388409726
$ws.Range("K17").Value = 1.03005213478444
$ws.Range("L17").Value = 1.024890698300706
$ws.Range("N17").Value = 1.012970743874219
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.020544222945653
$ws.Range("D18").Value = 1.026770942054959
$ws.Range("E18").Value = 1.021608204789362
$ws.Range("I18").Value = 1.030964463444248
$ws.Range("J18").Value = 1.026751077433335
$ws.Range("K18").Value = 1.030126928819723
$ws.Range("L18").Value = 1.024982400078018
$ws.Range("N18").Value = 1.012997138277852
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.020599200093904
$ws.Range("D19").Value = 1.02681103370723
$ws.Range("E19").Value = 1.021654177017656
$ws.Range("I19").Value = 1.030973593776413
$ws.Range("J19").Value = 1.02677824921246
$ws.Range("K19").Value = 1.030152428033672
$ws.Range("L19").Value = 1.025013674186894
$ws.Range("N19").Value = 1.013006137831535
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.020353385365508
$ws.Range("D20").Value = 1.026631777671978
$ws.Range("E20").Value = 1.021448658895468
$ws.Range("I20").Value = 1.030932690572477
$ws.Range("J20").Value = 1.02665673016711
$ws.Range("K20").Value = 1.030038375276824
$ws.Range("L20").Value = 1.024873833389681
$ws.Range("N20").Value = 1.012965888704258
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.019555210498447
$ws.Range("D21").Value = 1.026049767439358
$ws.Range("E21").Value = 1.02078192211061
$ws.Range("I21").Value = 1.030798480957844
$ws.Range("J21").Value = 1.026261657953079
$ws.Range("K21").Value = 1.029667334283365
$ws.Range("L21").Value = 1.024419632056684
$ws.Range("N21").Value = 1.012835019293771
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.019054138363309
$ws.Range("D22").Value = 1.025684434382329
$ws.Range("E22").Value = 1.020363819878631
$ws.Range("I22").Value = 1.030713159766647
$ws.Range("J22").Value = 1.026013266089096
$ws.Range("K22").Value = 1.029433866189773
$ws.Range("L22").Value = 1.024134398944343
$ws.Range("N22").Value = 1.012752727046908
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.019319692269349
$ws.Range("D23").Value = 1.025878046786356
$ws.Range("E23").Value = 1.020585358629439
$ws.Range("I23").Value = 1.030758478821999
$ws.Range("J23").Value = 1.026144942302652
$ws.Range("K23").Value = 1.029557648611628
$ws.Range("L23").Value = 1.024285573536175
$ws.Range("N23").Value = 1.012796352457979
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.020366779800013
$ws.Range("D24").Value = 1.026641545164373
$ws.Range("E24").Value = 1.021459855339376
$ws.Range("I24").Value = 1.030934924650942
$ws.Range("J24").Value = 1.026663353597562
$ws.Range("K24").Value = 1.030044592671369
$ws.Range("L24").Value = 1.02488145380393
$ws.Range("N24").Value = 1.012968082551033
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.021585155156572
$ws.Range("D25").Value = 1.027530083952755
$ws.Range("E25").Value = 1.022479383462418
$ws.Range("I25").Value = 1.03113556843479
$ws.Range("J25").Value = 1.027264918458275
$ws.Range("K25").Value = 1.03060883207594
$ws.Range("L25").Value = 1.025574372653441
$ws.Range("N25").Value = 1.013167308253146
